$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers in I1, J1 - copy formatting from the existing header H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New numeric data in columns I and J, rows 2-6 (unstyled, like column H)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7
